$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated K column (G) values: K count based on strikes instead of old Strike# metric.
$kValues = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 0
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    14 = 1
    15 = 1
    16 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
